$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Row 19 (S.No. 11) - "Java 8 features" topic, ht 409.5 (Excel's row-height cap)
# ---------------------------------------------------------------------------
$ws.Cells.Item(19,1).Value = 11

$ws.Cells.Item(19,2).Value = 44337
$ws.Cells.Item(19,2).NumberFormat = "d-mmm-yy"

$c19 = @"
Java 8 features- Lambda and Stream 
Java 8 features- Functional interfaces and default & static methods
"@
$ws.Cells.Item(19,3).Value = $c19
$ws.Cells.Item(19,3).WrapText = $true

$d19 = @"
https://www.tutorialspoint.com/java8/java8_streams.htm
https://www.geeksforgeeks.org/lambda-expressions-java-8/
https://www.programiz.com/java-programming/outputstream
https://www.programiz.com/java-programming/inputstream
https://www.javatpoint.com/java-lambda-expressions
https://beginnersbook.com/2017/10/java-8-interface-changes-default-method-and-static-method/
https://www.tutorialspoint.com/default-method-vs-static-method-in-an-interface-in-java#:~:text=Since%20Java8%20static%20methods%20and,can%20have%20a%20default%20implementation.&text=In%20short%2C%20you%20can%20access,objects%20of%20the%20implementing%20classes.
"@
$ws.Cells.Item(19,4).Value = $d19

$ws.Rows.Item(19).RowHeight = 409.5

# ---------------------------------------------------------------------------
# Row 20 - "JDBC introduction..." topic with hyperlink on D20
# ---------------------------------------------------------------------------
$ws.Cells.Item(20,2).Value = 44340
$ws.Cells.Item(20,2).NumberFormat = "d-mmm"

$ws.Cells.Item(20,3).Value = "JDBC introfuction and CRUD operation, Unit Testing- Junit test cases"
$ws.Cells.Item(20,3).WrapText = $true

$d20 = @"
https://www.codejava.net/java-se/jdbc/jdbc-tutorial-sql-insert-select-update-and-delete-examples
https://www.javatpoint.com/java-jdbc
https://www.javatpoint.com/junit-tutorial
"@
$ws.Cells.Item(20,4).Value = $d20
$ws.Hyperlinks.Add($ws.Cells.Item(20,4), "https://www.codejava.net/java-se/jdbc/jdbc-tutorial-sql-insert-select-update-and-delete-examples", [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, $d20)
$ws.Cells.Item(20,4).WrapText = $true

$ws.Rows.Item(20).RowHeight = 120

# ---------------------------------------------------------------------------
# Row 21 - "Java Design Patterns" topic
# ---------------------------------------------------------------------------
$ws.Cells.Item(21,2).Value = 44341
$ws.Cells.Item(21,2).NumberFormat = "d-mmm"

$ws.Cells.Item(21,3).Value = "Java Design Patterns"

$d21 = @"
https://www.freecodecamp.org/news/the-basic-design-patterns-all-developers-need-to-know/
https://www.javatpoint.com/design-patterns-in-java
"@
$ws.Cells.Item(21,4).Value = $d21

$ws.Rows.Item(21).RowHeight = 90

# ---------------------------------------------------------------------------
# Sheet view state - scroll near the bottom, leave behind the selection that
# was active when the workbook was last saved.
# ---------------------------------------------------------------------------
$ws.Activate()
$ws.Range("C34:C39").Select()
